$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: two new daily price entries are inserted at the top of the
# "Pepino ensalada" data block (row 265 onward), pushing the existing rows
# (265-337) down by two (-> 267-339), which also extends the used range.
$ws.Rows("265:266").Insert()

# Seed the two new rows from the (now shifted) rows directly below them -
# most columns (Mercado, Region, Categoria, Variedad, Calidad, Unidad,
# Origen, Kg/Unidades, Clasificacion) keep the same values as their
# neighbours; only a handful of cells need to change afterwards.
$ws.Rows("267").Copy()
$ws.Rows("265").PasteSpecial()
$ws.Rows("268").Copy()
$ws.Rows("266").PasteSpecial()

# Row 265 - new entry
$ws.Range("D265").Value = 44508
$ws.Range("J265").Value = 200
$ws.Range("K265").Value = 8000
$ws.Range("L265").Value = 8000
$ws.Range("M265").Value = 8000
$ws.Range("P265").Value = 133

# Row 266 - new entry
$ws.Range("D266").Value = 44508
$ws.Range("J266").Value = 300
$ws.Range("K266").Value = 10000
$ws.Range("L266").Value = 10000
$ws.Range("M266").Value = 10000
$ws.Range("O266").Value = "Región de O'Higgins"
$ws.Range("P266").Value = 167

Write-Output "done"
